$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(14, 1).Value = 9686.69
$ws.Cells.Item(14, 2).Value = 9792.4500000000007
$ws.Cells.Item(14, 3).Value = 281.06
$ws.Cells.Item(14, 4).Value = 284.08999999999997
$ws.Cells.Item(14, 5).Value = $true
$ws.Cells.Item(14, 6).Value = 1.08
$ws.Cells.Item(14, 7).Value = 42620.766377314816
$ws.Cells.Item(14, 8).Value = $false
